# Refitting NCDEs to individual patients (for manuscript figure)
# Adds a new "Label" column (H) with 0/1 values and refreshes the refitted
# D/E/F values for the first (Iterations=100) block of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell H1 = "Label" (copy formatting from G1 header) ---
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value() = "Label"

# --- Updated regression values for rows 2-11 (Iterations = 100 block) ---
$ws.Range("D2").Value() = 0.6079522988133341
$ws.Range("E2").Value() = 0.6079522988133341

$ws.Range("D3").Value() = 0.5370380911056414
$ws.Range("E3").Value() = 0.5370380911056414

$ws.Range("D4").Value() = 0.5210502797261171
$ws.Range("E4").Value() = 0.5210502797261171

$ws.Range("D5").Value() = 0.5990811618443884
$ws.Range("E5").Value() = 0.5990811618443884

$ws.Range("D6").Value() = 0.5940049022570391
$ws.Range("E6").Value() = 0.5940049022570391

$ws.Range("D7").Value() = 0.285210650652162
$ws.Range("E7").Value() = 0.714789349347838

$ws.Range("D8").Value() = 0.5920621593423706
$ws.Range("E8").Value() = 0.4079378406576294

$ws.Range("D9").Value() = 0.5691973535264673
$ws.Range("E9").Value() = 0.4308026464735327

$ws.Range("D10").Value() = 0.4930333532119519
$ws.Range("E10").Value() = 0.5069666467880481

$ws.Range("D11").Value() = 0.3331709628586969
$ws.Range("E11").Value() = 0.6668290371413031
$ws.Range("F11").Value() = 0.8406529426574707

# --- New "Label" column values (H2:H21), 1 = MDD patient, 0 = Control patient ---
$ws.Range("H2").Value() = 0
$ws.Range("H3").Value() = 0
$ws.Range("H4").Value() = 0
$ws.Range("H5").Value() = 0
$ws.Range("H6").Value() = 0
$ws.Range("H7").Value() = 1
$ws.Range("H8").Value() = 1
$ws.Range("H9").Value() = 1
$ws.Range("H10").Value() = 1
$ws.Range("H11").Value() = 1

$ws.Range("H12").Value() = 0
$ws.Range("H13").Value() = 0
$ws.Range("H14").Value() = 0
$ws.Range("H15").Value() = 0
$ws.Range("H16").Value() = 0
$ws.Range("H17").Value() = 1
$ws.Range("H18").Value() = 1
$ws.Range("H19").Value() = 1
$ws.Range("H20").Value() = 1
$ws.Range("H21").Value() = 1
